# Natmi LR-pairs table: recompute sending/target cluster breakdown
# to include the "ECs" cluster (per Dr Hou's advice), expanding the
# 2-sender x 3-target grid (rows 2-7) into a full 3-sender x 3-target
# grid (rows 2-10: ECs/FAPs/sCs x ECs/FAPs/sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col11a1"
$ws.Range("C2").Value = "Ddr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.029375
$ws.Range("H2").Value = 0.088125
$ws.Range("I2").Value = 0.01978727861543612
$ws.Range("J2").Value = 0.01978727861543612
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2763116666666667
$ws.Range("N2").Value = 0.828935
$ws.Range("O2").Value = 0.02083107478128044
$ws.Range("P2").Value = 0.02083107478128044
$ws.Range("Q2").Value = 0.008116655208333334
$ws.Range("R2").Value = 0.07304989687499999
$ws.Range("S2").Value = 0.0004121902805561811
$ws.Range("T2").Value = 0.0004121902805561811

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col11a1"
$ws.Range("C3").Value = "Ddr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.029375
$ws.Range("H3").Value = 0.088125
$ws.Range("I3").Value = 0.01978727861543612
$ws.Range("J3").Value = 0.01978727861543612
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.180798333333333
$ws.Range("N3").Value = 6.542395
$ws.Range("O3").Value = 0.1644098988384798
$ws.Range("P3").Value = 0.1644098988384798
$ws.Range("Q3").Value = 0.06406095104166666
$ws.Range("R3").Value = 0.5765485593749999
$ws.Range("S3").Value = 0.003253224475452666
$ws.Range("T3").Value = 0.003253224475452666

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col11a1"
$ws.Range("C4").Value = "Ddr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.029375
$ws.Range("H4").Value = 0.088125
$ws.Range("I4").Value = 0.01978727861543612
$ws.Range("J4").Value = 0.01978727861543612
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.807288
$ws.Range("N4").Value = 32.421864
$ws.Range("O4").Value = 0.8147590263802398
$ws.Range("P4").Value = 0.8147590263802398
$ws.Range("Q4").Value = 0.317464085
$ws.Range("R4").Value = 2.857176765
$ws.Range("S4").Value = 0.01612186385942727
$ws.Range("T4").Value = 0.01612186385942727

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Col11a1"
$ws.Range("C5").Value = "Ddr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.859432
$ws.Range("H5").Value = 2.578296
$ws.Range("I5").Value = 0.5789215467241361
$ws.Range("J5").Value = 0.5789215467241361
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2763116666666667
$ws.Range("N5").Value = 0.828935
$ws.Range("O5").Value = 0.02083107478128044
$ws.Range("P5").Value = 0.02083107478128044
$ws.Range("Q5").Value = 0.2374710883066667
$ws.Range("R5").Value = 2.13723979476
$ws.Range("S5").Value = 0.01205955803230502
$ws.Range("T5").Value = 0.01205955803230502

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col11a1"
$ws.Range("C6").Value = "Ddr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.859432
$ws.Range("H6").Value = 2.578296
$ws.Range("I6").Value = 0.5789215467241361
$ws.Range("J6").Value = 0.5789215467241361
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.180798333333333
$ws.Range("N6").Value = 6.542395
$ws.Range("O6").Value = 0.1644098988384798
$ws.Range("P6").Value = 0.1644098988384798
$ws.Range("Q6").Value = 1.874247873213333
$ws.Range("R6").Value = 16.86823085892
$ws.Range("S6").Value = 0.09518043293233144
$ws.Range("T6").Value = 0.09518043293233144

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col11a1"
$ws.Range("C7").Value = "Ddr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.859432
$ws.Range("H7").Value = 2.578296
$ws.Range("I7").Value = 0.5789215467241361
$ws.Range("J7").Value = 0.5789215467241361
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 10.807288
$ws.Range("N7").Value = 32.421864
$ws.Range("O7").Value = 0.8147590263802398
$ws.Range("P7").Value = 0.8147590263802398
$ws.Range("Q7").Value = 9.288129140416
$ws.Range("R7").Value = 83.59316226374399
$ws.Range("S7").Value = 0.4716815557594996
$ws.Range("T7").Value = 0.4716815557594996

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Col11a1"
$ws.Range("C8").Value = "Ddr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.5957326666666666
$ws.Range("H8").Value = 1.787198
$ws.Range("I8").Value = 0.4012911746604278
$ws.Range("J8").Value = 0.4012911746604279
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2763116666666667
$ws.Range("N8").Value = 0.828935
$ws.Range("O8").Value = 0.02083107478128044
$ws.Range("P8").Value = 0.02083107478128044
$ws.Range("Q8").Value = 0.1646078860144444
$ws.Range("R8").Value = 1.48147097413
$ws.Range("S8").Value = 0.008359326468419241
$ws.Range("T8").Value = 0.008359326468419243

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Col11a1"
$ws.Range("C9").Value = "Ddr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.5957326666666666
$ws.Range("H9").Value = 1.787198
$ws.Range("I9").Value = 0.4012911746604278
$ws.Range("J9").Value = 0.4012911746604279
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.180798333333333
$ws.Range("N9").Value = 6.542395
$ws.Range("O9").Value = 0.1644098988384798
$ws.Range("P9").Value = 0.1644098988384798
$ws.Range("Q9").Value = 1.299172806578889
$ws.Range("R9").Value = 11.69255525921
$ws.Range("S9").Value = 0.06597624143069565
$ws.Range("T9").Value = 0.06597624143069565

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col11a1"
$ws.Range("C10").Value = "Ddr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5957326666666666
$ws.Range("H10").Value = 1.787198
$ws.Range("I10").Value = 0.4012911746604278
$ws.Range("J10").Value = 0.4012911746604279
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 10.807288
$ws.Range("N10").Value = 32.421864
$ws.Range("O10").Value = 0.8147590263802398
$ws.Range("P10").Value = 0.8147590263802398
$ws.Range("Q10").Value = 6.438254499674666
$ws.Range("R10").Value = 57.944290497072
$ws.Range("S10").Value = 0.3269556067613129
$ws.Range("T10").Value = 0.326955606761313
